# posts.xlsx — remove the post row about the tiger-cub/napping-mother
# photo (old row 652). Excel's Rows().Delete() removes the entire row
# and shifts every row beneath it up by one, which is exactly what the
# diff shows: row 652 disappears and rows 653..737 become 652..736
# (their A/B/C content is otherwise untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(652).Delete()
